$d = $word.ActiveDocument

# ============================================================================
# Edit 2 (do first, so it doesn't shift offsets used by Edit 1 below):
#   "{statusTanah}" + " - "  ->  single run "{statusTanah} - "
#   The following "Girik/Persil/.../lainnya*" run must stay a separate run.
# ============================================================================

$rng2 = $d.Content
$rng2.Find.Execute("{statusTanah}")
$start2 = $rng2.Start
$end2 = $rng2.End
$sepStart = $end2         # start of the " - " run
$sepEnd = $end2 + 3       # end of the " - " run
$afterStart = $sepEnd     # start of the following "Girik/..." run

# Isolate the following run (the one that would otherwise also get merged in)
# by giving it a momentarily-different character format, then restoring it;
# formatting-only changes do not trigger the engine's run-coalescing pass.
$afterLen = 46  # length of "Girik/Persil/AJB/SHM/SHGB/Tanah Adat/lainnya*"
$afterRng = $d.Range($afterStart, $afterStart + $afterLen)
$afterRng.Bold = 1

# Delete the " - " run and reinsert identical text. Because this is a genuine
# insert/delete mutation (not a same-value no-op), the engine performs its
# run-coalescing pass and merges the preceding "{statusTanah}" run with the
# newly (re)inserted text into one run, but stops there since the following
# run's format currently differs.
$sepRng = $d.Range($sepStart, $sepEnd)
$sepRng.Delete()
$sepIns = $d.Range($sepStart, $sepStart)
$sepIns.InsertBefore(" - ")

# Restore the following run's formatting (formatting-only change -> no merge)
$afterRng2 = $d.Range($afterStart, $afterStart + $afterLen)
$afterRng2.Bold = 0

# ============================================================================
# Edit 1:
#   "{tempatLahir}, {ttl}"  ->  5 runs:
#     "{tempatLahir" | "Pembuat" | "}, {ttl" | "Pembuat" | "}"
# ============================================================================

$rng1 = $d.Content
$rng1.Find.Execute("{tempatLahir}, {ttl}")
$s1 = $rng1.Start

# Offsets (relative to $s1) of the two "}" characters in "{tempatLahir}, {ttl}"
$firstBracePos = $s1 + 12    # right after "{tempatLahir"
$secondBracePos = $s1 + 19   # right after "{tempatLahir}, {ttl"

# Insert the first "Pembuat" (before the first "}") and mark it bold so the
# later edit (second insertion) cannot silently re-merge it with its
# right-hand neighbour; we restore formatting only once both insertions
# are complete.
$ins1 = $d.Range($firstBracePos, $firstBracePos)
$ins1.InsertBefore("Pembuat")
$run1 = $d.Range($firstBracePos, $firstBracePos + 7)
$run1.Bold = 1

# Insert the second "Pembuat" (before the second "}"); the text shifted
# right by 7 characters ("Pembuat") because of the first insertion.
$secondBracePos2 = $secondBracePos + 7
$ins2 = $d.Range($secondBracePos2, $secondBracePos2)
$ins2.InsertBefore("Pembuat")
$run2 = $d.Range($secondBracePos2, $secondBracePos2 + 7)
$run2.Bold = 1

# Restore both inserted runs to normal (non-bold) formatting. These are
# formatting-only changes, so they do not trigger another run-coalescing
# pass and the five runs created above remain distinct.
$d.Range($firstBracePos, $firstBracePos + 7).Bold = 0
$d.Range($secondBracePos2, $secondBracePos2 + 7).Bold = 0

Write-Output "Done"
